$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.276.73"
$ws.Cells.Item(2, 5).Value = "  +0.45%  "
$ws.Cells.Item(3, 4).Value = "3.113.87"
$ws.Cells.Item(3, 5).Value = "  +0.03%  "
$ws.Cells.Item(5, 4).Value = "'579.76"
$ws.Cells.Item(5, 5).Value = "  +0.04%  "
$ws.Cells.Item(6, 4).Value = "'174.42"
$ws.Cells.Item(6, 5).Value = "  +0.92%  "
$ws.Cells.Item(7, 5).Value = "  -0.07%  "
$ws.Cells.Item(8, 4).Value = "'0.520"
$ws.Cells.Item(8, 5).Value = "  -0.65%  "
$ws.Cells.Item(9, 5).Value = "  +1.44%  "
$ws.Cells.Item(10, 5).Value = "  -1.06%  "
$ws.Cells.Item(11, 5).Value = "  -0.85%  "
$ws.Cells.Item(12, 5).Value = "  -0.73%  "
$ws.Cells.Item(13, 4).Value = "'36.86"
$ws.Cells.Item(13, 5).Value = "  -1.21%  "
$ws.Cells.Item(14, 5).Value = "  -1.80%  "
$ws.Cells.Item(15, 4).Value = "3.629.55"
$ws.Cells.Item(15, 5).Value = "  +0.05%  "
$ws.Cells.Item(16, 4).Value = "67.209.60"
$ws.Cells.Item(16, 5).Value = "  +0.28%  "
$ws.Cells.Item(17, 5).Value = "  -1.45%  "
$ws.Cells.Item(18, 2).Value = "WrappedEther"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(18, 4).Value = "3.114.81"
$ws.Cells.Item(18, 5).Value = "  -0.01%  "
$ws.Cells.Item(19, 2).Value = "Chainlink"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(19, 4).Value = "'16.63"
$ws.Cells.Item(19, 5).Value = "  +1.84%  "
$ws.Cells.Item(20, 4).Value = "'492.80"
$ws.Cells.Item(20, 5).Value = "  +1.19%  "
$ws.Cells.Item(21, 4).Value = "'0.704"
$ws.Cells.Item(21, 5).Value = "  -2.26%  "
$ws.Cells.Item(22, 5).Value = "  +3.61%  "
$ws.Cells.Item(23, 4).Value = "'84.01"
$ws.Cells.Item(23, 5).Value = "  -0.71%  "
$ws.Cells.Item(24, 4).Value = "'13.13"
$ws.Cells.Item(24, 5).Value = "  -1.60%  "
$ws.Cells.Item(25, 5).Value = "  -3.04%  "
$ws.Cells.Item(26, 4).Value = "'10.57"
$ws.Cells.Item(26, 5).Value = "  +4.82%  "
$ws.Cells.Item(28, 4).Value = "'7.93"
$ws.Cells.Item(28, 5).Value = "  -1.58%  "
$ws.Cells.Item(29, 5).Value = "  -2.79%  "
$ws.Cells.Item(30, 5).Value = "  -0.68%  "
$ws.Cells.Item(31, 4).Value = "'28.42"
$ws.Cells.Item(31, 5).Value = "  -2.07%  "
$ws.Cells.Item(32, 5).Value = "  -0.98%  "
$ws.Cells.Item(34, 5).Value = "  -0.09%  "
$ws.Cells.Item(35, 4).Value = "'5.85"
$ws.Cells.Item(35, 5).Value = "  -1.67%  "
$ws.Cells.Item(36, 5).Value = "  -1.55%  "
$ws.Cells.Item(37, 4).Value = "'47.22"
$ws.Cells.Item(37, 5).Value = "  -0.38%  "
$ws.Cells.Item(38, 5).Value = "  -3.42%  "
$ws.Cells.Item(39, 5).Value = "  -2.23%  "
$ws.Cells.Item(40, 5).Value = "  +0.81%  "
$ws.Cells.Item(41, 5).Value = "  -2.48%  "
$ws.Cells.Item(42, 4).Value = "'390.21"
$ws.Cells.Item(42, 5).Value = "  +1.13%  "
$ws.Cells.Item(43, 4).Value = "2.801.19"
$ws.Cells.Item(43, 5).Value = "  -1.51%  "
$ws.Cells.Item(44, 4).Value = "'2.56"
$ws.Cells.Item(44, 5).Value = "  -8.46%  "
$ws.Cells.Item(45, 5).Value = "  -2.15%  "
$ws.Cells.Item(46, 4).Value = "'134.96"
$ws.Cells.Item(46, 5).Value = "  -1.03%  "
$ws.Cells.Item(47, 5).Value = "  +0.02%  "
$ws.Cells.Item(48, 5).Value = "  +0.42%  "
$ws.Cells.Item(49, 5).Value = "  -1.16%  "
$ws.Cells.Item(50, 5).Value = "  -1.17%  "
$ws.Cells.Item(51, 4).Value = "'6.72"
$ws.Cells.Item(51, 5).Value = "  -2.65%  "
